$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Fix D column values (8 -> 6) for morici rows 187,190,193,196,199,202,205,208
$fixRows = @(187,190,193,196,199,202,205,208)
foreach ($r in $fixRows) {
    $ws.Cells.Item($r, 4).Value = 6
}

# Append new Morici burn-treatment rows (210-233)
# Columns: A value, B variable, C units, D years_post, E treatment,
# F study, G forest_type, H region, I burn_season, J thin_type, K notes
$newRows = @(
    @("1998.09, 10.02","all_woody","Mg/ha",0,"burn","morici","ponderosa","interior_pnw","fall","NA","NA"),
    @("2006.64, 4.14","all_woody","Mg/ha",4,"burn","morici","ponderosa","interior_pnw","fall","NA","NA"),
    @("2015.27, 8.82","all_woody","Mg/ha",15,"burn","morici","ponderosa","interior_pnw","fall","NA","NA"),
    @("1998.07, 7.52","cwd","Mg/ha",0,"burn","morici","ponderosa","interior_pnw","fall","NA","NA"),
    @("2006.43, 2.45","cwd","Mg/ha",4,"burn","morici","ponderosa","interior_pnw","fall","NA","NA"),
    @("2015.25, 6.75","cwd","Mg/ha",15,"burn","morici","ponderosa","interior_pnw","fall","NA","NA"),
    @("1997.86, 5.07","cwd_rotten","Mg/ha",0,"burn","morici","ponderosa","interior_pnw","fall","NA","NA"),
    @("2006.42, 0.44","cwd_rotten","Mg/ha",4,"burn","morici","ponderosa","interior_pnw","fall","NA","NA"),
    @("2015.04, 4.08","cwd_rotten","Mg/ha",15,"burn","morici","ponderosa","interior_pnw","fall","NA","NA"),
    @("1998.03, 2.61","fwd","Mg/ha",0,"burn","morici","ponderosa","interior_pnw","fall","NA","NA"),
    @("2006.43, 1.63","fwd","Mg/ha",4,"burn","morici","ponderosa","interior_pnw","fall","NA","NA"),
    @("2015.02, 2.07","fwd","Mg/ha",15,"burn","morici","ponderosa","interior_pnw","fall","NA","NA"),
    @("1997.84, 2.40","cwd_sound","Mg/ha",0,"burn","morici","ponderosa","interior_pnw","fall","NA","NA"),
    @("2006.43, 2.07","cwd_sound","Mg/ha",4,"burn","morici","ponderosa","interior_pnw","fall","NA","NA"),
    @("2015.03, 2.61","cwd_sound","Mg/ha",15,"burn","morici","ponderosa","interior_pnw","fall","NA","NA"),
    @("1997.83, 1.91","hundred_hour","Mg/ha",0,"burn","morici","ponderosa","interior_pnw","fall","NA","NA"),
    @("2006.42, 1.09","hundred_hour","Mg/ha",4,"burn","morici","ponderosa","interior_pnw","fall","NA","NA"),
    @("2015.02, 1.20","hundred_hour","Mg/ha",15,"burn","morici","ponderosa","interior_pnw","fall","NA","NA"),
    @("1997.82, 0.60","ten_hour","Mg/ha",0,"burn","morici","ponderosa","interior_pnw","fall","NA","NA"),
    @("2006.42, 0.38","ten_hour","Mg/ha",4,"burn","morici","ponderosa","interior_pnw","fall","NA","NA"),
    @("2015.01, 0.76","ten_hour","Mg/ha",15,"burn","morici","ponderosa","interior_pnw","fall","NA","NA"),
    @("1997.82, 0.05","one_hour","Mg/ha",0,"burn","morici","ponderosa","interior_pnw","fall","NA","NA"),
    @("2006.41, 0.00","one_hour","Mg/ha",4,"burn","morici","ponderosa","interior_pnw","fall","NA","NA"),
    @("2015.01, 0.11","one_hour","Mg/ha",15,"burn","morici","ponderosa","interior_pnw","fall","NA","NA")
)

$startRow = 210
$numCols = 11
# Write column-by-column so shared strings are interned in the same order
# as the original authoring (all col-A coordinate strings before "burn").
for ($c = 0; $c -lt $numCols; $c++) {
    for ($n = 0; $n -lt $newRows.Count; $n++) {
        $row = $startRow + $n
        $ws.Cells.Item($row, $c + 1).Value = $newRows[$n][$c]
    }
}

# Update sheet view: freeze top row, set top-left cell and selection
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.ScrollRow = 202
$ws.Range("G229").Select()
